$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '69.806.31'
$ws.Range("E2").Value = '  -1.72%  '
$ws.Range("D3").Value = '3.775.19'
$ws.Range("E3").Value = '  +2.51%  '
$ws.Range("E4").Value = '  +0.08%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '619.65'
$ws.Range("E5").Value = '  +3.46%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '177.05'
$ws.Range("E6").Value = '  -4.15%  '
$ws.Range("D7").Value = '3.772.57'
$ws.Range("E7").Value = '  +2.52%  '
$ws.Range("E8").Value = '  +0.05%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.534'
$ws.Range("E9").Value = '  -0.54%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.169'
$ws.Range("E10").Value = '  +3.72%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '6.27'
$ws.Range("E11").Value = '  -5.80%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.490'
$ws.Range("E12").Value = '  -2.09%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '40.71'
$ws.Range("E13").Value = '  +1.31%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.0000259'
$ws.Range("E14").Value = '  +1.90%  '
$ws.Range("D15").Value = '4.404.51'
$ws.Range("E15").Value = '  +2.58%  '
$ws.Range("D16").Value = '3.777.49'
$ws.Range("D17").Value = '69.860.14'
$ws.Range("E17").Value = '  -1.76%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.123'
$ws.Range("E18").Value = '  +0.29%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '7.56'
$ws.Range("E19").Value = '  +0.69%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '16.76'
$ws.Range("E20").Value = '  -1.53%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '507.44'
$ws.Range("E21").Value = '  -1.36%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '9.60'
$ws.Range("E22").Value = '  +4.23%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.724'
$ws.Range("E23").Value = '  -2.52%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.51'
$ws.Range("E24").Value = '  +1.48%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '87.05'
$ws.Range("E25").Value = '  -0.58%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '13.15'
$ws.Range("E26").Value = '  -2.77%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '11.04'
$ws.Range("E27").Value = '  -0.42%  '
$ws.Range("E28").Value = '  +24.75%  '
$ws.Range("E29").Value = '  +0.09%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '2.48'
$ws.Range("E30").Value = '  -1.58%  '
$ws.Range("E31").Value = '  +4.64%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '7.80'
$ws.Range("E32").Value = '  -4.54%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '31.07'
$ws.Range("E33").Value = '  -1.80%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.114'
$ws.Range("E34").Value = '  -2.07%  '
$ws.Range("E35").Value = '  +0.10%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.07'
$ws.Range("E36").Value = '  +5.45%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '6.17'
$ws.Range("E37").Value = '  +0.77%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.333'
$ws.Range("E38").Value = '  -3.52%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.132'
$ws.Range("E39").Value = '  +3.09%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.11'
$ws.Range("E40").Value = '  -2.48%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '50.89'
$ws.Range("E41").Value = '  -0.22%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '45.47'
$ws.Range("E42").Value = '  -0.10%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '422.67'
$ws.Range("E43").Value = '  +2.15%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '8.69'
$ws.Range("E44").Value = '  -1.71%  '
$ws.Range("D45").Value = '3.034.03'
$ws.Range("E45").Value = '  -3.94%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.79'
$ws.Range("E46").Value = '  +0.01%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0361'
$ws.Range("E47").Value = '  -2.03%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '27.24'
$ws.Range("E48").Value = '  -3.72%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '138.61'
$ws.Range("E49").Value = '  +0.62%  '
$ws.Range("E50").Value = '  -0.04%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.47'
$ws.Range("E51").Value = '  +0.18%  '
